$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 278
$wsOff.Range("C3").Value = 197
$wsOff.Range("D3").Value = 74
$wsOff.Range("E3").Value = 37

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 326
$wsDef.Range("C3").Value = 240
$wsDef.Range("D3").Value = 77
$wsDef.Range("E3").Value = 37
$wsDef.Range("G3").Value = 3
